$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# Update Status text for both language sheets (shared string "Ready for handoff" -> "Handoff transform failed")
$wsZh.Range("B2").Value = "Handoff transform failed"
$wsDe.Range("B2").Value = "Handoff transform failed"

foreach ($ws in @($wsZh, $wsDe)) {
    # Remove the "Latest Handoff File" hyperlink + value in C2 (handoff transform failed, so no file)
    $ws.Range("C2").Hyperlinks.Delete()
    $ws.Range("C2").ClearContents()

    # Reset "Latest Handoff Datetime" to epoch value
    $ws.Range("D2").Value = "0001-01-01 00:00:00"

    # Update "Handoff Reason" to Ignored
    $ws.Range("H2").Value = "Ignored"
}

$wb.Save()
